# 自动更新Excel文件 - 每日剩余天数递减，到期后重置为总天数并更新开始时间为今天
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

$newStartDate = 20260127

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2     # D 列：总天数
    $remain = $ws.Cells.Item($r, 5).Value2    # E 列：剩余天数
    $startDate = $ws.Cells.Item($r, 6).Value2 # F 列：开始时间 (YYYYMMDD)

    if ($total -eq $null -or $remain -eq $null) {
        continue
    }

    # 校验开始时间是否为合法的8位日期格式，非法（如格式损坏）则跳过该行
    $dateStr = [string]$startDate
    if ($dateStr.Length -ne 8) {
        continue
    }

    if ($remain -le 1) {
        # 剩余天数到期，重置为总天数，并将开始时间更新为今天
        $ws.Cells.Item($r, 5).Value = $total
        $ws.Cells.Item($r, 6).Value = $newStartDate
    } else {
        # 剩余天数递减1
        $ws.Cells.Item($r, 5).Value = $remain - 1
    }
}
